# ============================================================================
# Adds a new daily-snapshot sheet ('동원산업') to the workbook and appends the
# latest five trading days (2025-11-11 .. 2025-11-14, i.e. serials 45971-45975)
# to the three existing snapshot sheets, replacing each sheet's previously-zero
# placeholder value for the most recent existing day (serial 45968 / row 105).
# ============================================================================

$wb = $excel.ActiveWorkbook

# xlPasteFormats: used after writing raw values so the new cells pick up the
# workbook's existing numeric date-format style (no new style entries created).
$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# Part 1: update the three existing sheets (한화솔루션, 아난티, 대아티아이).
# Each sheet gets the same five new dates in rows 106-110, plus an update to
# the existing row 105's remn_amt (was a 0 placeholder).
# ---------------------------------------------------------------------------
$sheetUpdates = @(
    @(1, @(@(105,45968,439225), @(106,45971,445963), @(107,45972,466990), @(108,45973,478631), @(109,45974,489509), @(110,45975,470563))),
    @(2, @(@(105,45968,57483), @(106,45971,57287), @(107,45972,57076), @(108,45973,56788), @(109,45974,57444), @(110,45975,57084))),
    @(3, @(@(105,45968,13097), @(106,45971,13412), @(107,45972,13245), @(108,45973,13559), @(109,45974,13611), @(110,45975,13272)))
)

foreach ($entry in $sheetUpdates) {
    $sheetIndex = $entry[0]
    $newRows = $entry[1]
    $ws = $wb.Worksheets.Item($sheetIndex)

    foreach ($row in $newRows) {
        $r = $row[0]; $d = $row[1]; $v = $row[2]
        $ws.Cells.Item($r, 1).Value = $d
        $ws.Cells.Item($r, 2).Value = $v
    }

    # Row 105 (serial 45968) already carries the correct date-format style;
    # copy it onto the brand-new rows 106-110 so column A keeps its format.
    $ws.Range("A105").Copy()
    $ws.Range("A106:A110").PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------------
# Part 2: append the new '동원산업' worksheet as the last tab in the workbook
# (sheetId 4 / rId4), mirroring the date / remn_amt layout of the other sheets.
# ---------------------------------------------------------------------------
$wsTemplate = $wb.Worksheets.Item(1)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $lastSheet)
$ws4.Name = "동원산업"

$ws4.Range("A1").Value = "date"
$ws4.Range("B1").Value = "remn_amt"

# Match the bold/centered/bordered header style used by the other sheets.
$wsTemplate.Range("A1:B1").Copy()
$ws4.Range("A1:B1").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Part 3: populate '동원산업' data rows 2-101 (2025-06-20 .. 2025-11-14).
# ---------------------------------------------------------------------------
$sheet4Data = @(
    @(45828,31227),
    @(45831,30323),
    @(45832,29731),
    @(45833,30004),
    @(45834,30827),
    @(45835,30193),
    @(45838,31251),
    @(45839,33704),
    @(45840,35018),
    @(45841,31947),
    @(45842,30359),
    @(45845,30922),
    @(45846,30505),
    @(45847,33411),
    @(45848,31345),
    @(45849,30384),
    @(45852,30756),
    @(45853,30616),
    @(45854,29876),
    @(45855,29856),
    @(45856,29916),
    @(45859,29227),
    @(45860,29699),
    @(45861,28705),
    @(45862,27716),
    @(45863,27462),
    @(45866,26552),
    @(45867,26983),
    @(45868,27918),
    @(45869,16336),
    @(45870,10094),
    @(45873,11238),
    @(45874,15081),
    @(45875,15116),
    @(45876,15186),
    @(45877,15961),
    @(45880,15796),
    @(45881,15709),
    @(45882,16935),
    @(45883,19440),
    @(45887,18877),
    @(45888,18125),
    @(45889,18024),
    @(45890,18129),
    @(45891,17713),
    @(45894,18291),
    @(45895,18612),
    @(45896,18495),
    @(45897,20319),
    @(45898,20470),
    @(45901,20015),
    @(45902,13916),
    @(45903,14093),
    @(45904,15089),
    @(45905,14940),
    @(45908,14697),
    @(45909,15142),
    @(45910,16580),
    @(45911,16620),
    @(45912,16769),
    @(45915,17059),
    @(45916,17030),
    @(45917,16322),
    @(45918,16339),
    @(45919,16358),
    @(45922,17528),
    @(45923,13563),
    @(45924,12940),
    @(45925,12791),
    @(45926,12124),
    @(45929,12059),
    @(45930,12138),
    @(45931,12098),
    @(45932,11712),
    @(45940,13236),
    @(45943,11492),
    @(45944,11086),
    @(45945,11462),
    @(45946,11668),
    @(45947,11541),
    @(45950,11665),
    @(45951,11417),
    @(45952,12132),
    @(45953,12464),
    @(45954,13234),
    @(45957,14048),
    @(45958,13430),
    @(45959,13163),
    @(45960,12793),
    @(45961,13051),
    @(45964,12922),
    @(45965,12779),
    @(45966,12572),
    @(45967,13307),
    @(45968,12926),
    @(45971,13457),
    @(45972,15039),
    @(45973,15275),
    @(45974,16857),
    @(45975,14718)
)

$r = 2
foreach ($row in $sheet4Data) {
    $ws4.Cells.Item($r, 1).Value = $row[0]
    $ws4.Cells.Item($r, 2).Value = $row[1]
    $r++
}

# Column A (dates) picks up the same numeric date-format style as the other sheets.
$wsTemplate.Range("A2").Copy()
$ws4.Range("A2:A101").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
